$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, 8).Value = 84
$ws.Cells.Item(4, 10).Value = 99.14286
$ws.Cells.Item(4, 12).Value = 99.14286
$ws.Cells.Item(4, 14).Value = -327.14286
$ws.Cells.Item(39, 8).Value = 467.16666
$ws.Cells.Item(39, 9).Value = 78.5
$ws.Cells.Item(39, 11).Value = 235.5
$ws.Cells.Item(39, 13).Value = 60.5
$ws.Cells.Item(40, 8).Value = 5176.0884
$ws.Cells.Item(40, 9).Value = 1829
$ws.Cells.Item(40, 10).Value = 5500
$ws.Cells.Item(40, 11).Value = 1829
$ws.Cells.Item(40, 12).Value = 5500
$ws.Cells.Item(40, 13).Value = -1654
$ws.Cells.Item(40, 14).Value = -5850
$ws.Cells.Item(51, 8).Value = 2995.5
$ws.Cells.Item(51, 9).Value = 2995.2632
$ws.Cells.Item(51, 11).Value = 2995.2632
$ws.Cells.Item(51, 13).Value = -2511.2632
$ws.Cells.Item(86, 8).Value = 6785.25
$ws.Cells.Item(86, 9).Value = 6021.1113
$ws.Cells.Item(86, 10).Value = 8160.7
$ws.Cells.Item(86, 11).Value = 6021.1113
$ws.Cells.Item(86, 12).Value = 8160.7
$ws.Cells.Item(86, 13).Value = -4898.1113
$ws.Cells.Item(86, 14).Value = -10406.7
$ws.Cells.Item(89, 8).Value = 6785.25
$ws.Cells.Item(89, 9).Value = 6021.1113
$ws.Cells.Item(89, 10).Value = 8160.7
$ws.Cells.Item(89, 11).Value = 30105.5565
$ws.Cells.Item(89, 12).Value = 40803.5
$ws.Cells.Item(89, 13).Value = -24489.5565
$ws.Cells.Item(89, 14).Value = -52035.5
$ws.Cells.Item(121, 8).Value = 4950
$ws.Cells.Item(121, 10).Value = 4950
$ws.Cells.Item(121, 12).Value = 14850
$ws.Cells.Item(121, 14).Value = -18344
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(23, 8).Value = 6400
$ws.Cells.Item(23, 10).Value = 7500
$ws.Cells.Item(23, 12).Value = 7500
$ws.Cells.Item(23, 14).Value = -8018
$ws.Cells.Item(32, 8).Value = 1428139.9
$ws.Cells.Item(32, 9).Value = 1635113.5
$ws.Cells.Item(32, 10).Value = 20719.2
$ws.Cells.Item(32, 11).Value = 1635113.5
$ws.Cells.Item(32, 12).Value = 20719.2
$ws.Cells.Item(32, 13).Value = -1634826.5
$ws.Cells.Item(32, 14).Value = -21293.2
$ws.Cells.Item(61, 8).Value = 3032876.2
$ws.Cells.Item(61, 9).Value = 2441.6667
$ws.Cells.Item(61, 10).Value = 16669832
$ws.Cells.Item(61, 11).Value = 2441.6667
$ws.Cells.Item(61, 12).Value = 16669832
$ws.Cells.Item(61, 13).Value = -2229.6667
$ws.Cells.Item(61, 14).Value = -16670256
$ws.Cells.Item(122, 8).Value = 2059.9443
$ws.Cells.Item(122, 9).Value = 1739.3077
$ws.Cells.Item(122, 10).Value = 2893.6
$ws.Cells.Item(122, 11).Value = 5217.9231
$ws.Cells.Item(122, 12).Value = 8680.799999999999
$ws.Cells.Item(122, 13).Value = -2767.9231
$ws.Cells.Item(122, 14).Value = -13580.8
$ws.Cells.Item(132, 8).Value = 3870.173
$ws.Cells.Item(132, 9).Value = 2413.0334
$ws.Cells.Item(132, 10).Value = 5857.1816
$ws.Cells.Item(132, 11).Value = 7239.100199999999
$ws.Cells.Item(132, 12).Value = 17571.5448
$ws.Cells.Item(132, 13).Value = -4709.100199999999
$ws.Cells.Item(132, 14).Value = -22631.5448
$ws.Cells.Item(136, 8).Value = 3032876.2
$ws.Cells.Item(136, 9).Value = 2441.6667
$ws.Cells.Item(136, 10).Value = 16669832
$ws.Cells.Item(136, 11).Value = 7325.000100000001
$ws.Cells.Item(136, 12).Value = 50009496
$ws.Cells.Item(136, 13).Value = -4775.000100000001
$ws.Cells.Item(136, 14).Value = -50014596
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(35, 8).Value = 31416.334
$ws.Cells.Item(35, 10).Value = 31416.334
$ws.Cells.Item(35, 12).Value = 31416.334
$ws.Cells.Item(35, 14).Value = -32036.334
$ws.Cells.Item(94, 8).Value = 2817.7222
$ws.Cells.Item(94, 9).Value = 3222.2666
$ws.Cells.Item(94, 11).Value = 3222.2666
$ws.Cells.Item(94, 13).Value = -2771.2666
$ws.Cells.Item(134, 8).Value = 2528870
$ws.Cells.Item(134, 9).Value = 3730.4375
$ws.Cells.Item(134, 10).Value = 83333336
$ws.Cells.Item(134, 11).Value = 11191.3125
$ws.Cells.Item(134, 12).Value = 250000008
$ws.Cells.Item(134, 13).Value = -8656.3125
$ws.Cells.Item(134, 14).Value = -250005078
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(10, 8).Value = 828.5
$ws.Cells.Item(10, 9).Value = 828.5
$ws.Cells.Item(10, 10).Value = 0
$ws.Cells.Item(10, 11).Value = 828.5
$ws.Cells.Item(10, 12).Value = 0
$ws.Cells.Item(10, 13).Value = $null
$ws.Cells.Item(10, 14).Value = -689.5
$ws.Cells.Item(58, 8).Value = 6097134.5
$ws.Cells.Item(58, 9).Value = 11906575
$ws.Cells.Item(58, 10).Value = 2560953
$ws.Cells.Item(58, 11).Value = 11906575
$ws.Cells.Item(58, 12).Value = 2560953
$ws.Cells.Item(58, 13).Value = -11906372
$ws.Cells.Item(58, 14).Value = -2561359
$ws.Cells.Item(132, 8).Value = 5491
$ws.Cells.Item(132, 9).Value = 5762.5
$ws.Cells.Item(132, 11).Value = 17287.5
$ws.Cells.Item(132, 13).Value = -14757.5
$ws.Cells.Item(134, 8).Value = 2532.4119
$ws.Cells.Item(134, 9).Value = 2571.9375
$ws.Cells.Item(134, 10).Value = 1900
$ws.Cells.Item(134, 11).Value = 7715.8125
$ws.Cells.Item(134, 12).Value = 5700
$ws.Cells.Item(134, 13).Value = -5180.8125
$ws.Cells.Item(134, 14).Value = -10770
$ws.Cells.Item(136, 8).Value = 6097134.5
$ws.Cells.Item(136, 9).Value = 11906575
$ws.Cells.Item(136, 10).Value = 2560953
$ws.Cells.Item(136, 11).Value = 35719725
$ws.Cells.Item(136, 12).Value = 7682859
$ws.Cells.Item(136, 13).Value = -35717175
$ws.Cells.Item(136, 14).Value = -7687959
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(57, 8).Value = 4002.5
$ws.Cells.Item(57, 9).Value = 4002.5
$ws.Cells.Item(57, 11).Value = 12007.5
$ws.Cells.Item(57, 13).Value = -11448.5
$ws.Cells.Item(68, 8).Value = 6698.1787
$ws.Cells.Item(68, 10).Value = 6963.423
$ws.Cells.Item(68, 12).Value = 20890.269
$ws.Cells.Item(68, 14).Value = -22512.269
$ws.Cells.Item(71, 8).Value = 6698.1787
$ws.Cells.Item(71, 10).Value = 6963.423
$ws.Cells.Item(71, 12).Value = 62670.807
$ws.Cells.Item(71, 14).Value = -70782.807
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 54750
$ws.Cells.Item(70, 10).Value = 20000
$ws.Cells.Item(70, 12).Value = 20000
$ws.Cells.Item(70, 14).Value = -20540
$ws.Cells.Item(73, 8).Value = 54750
$ws.Cells.Item(73, 10).Value = 20000
$ws.Cells.Item(73, 12).Value = 20000
$ws.Cells.Item(73, 14).Value = -21872
$ws.Cells.Item(102, 8).Value = 1530.3
$ws.Cells.Item(102, 9).Value = 911.13336
$ws.Cells.Item(102, 10).Value = 3387.8
$ws.Cells.Item(102, 11).Value = 911.13336
$ws.Cells.Item(102, 12).Value = 3387.8
$ws.Cells.Item(102, 13).Value = 710.86664
$ws.Cells.Item(102, 14).Value = -6631.8
$ws.Cells.Item(122, 8).Value = 39891.703
$ws.Cells.Item(122, 9).Value = 46181.566
$ws.Cells.Item(122, 11).Value = 138544.698
$ws.Cells.Item(122, 13).Value = -136094.698
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 7408.875
$ws.Cells.Item(40, 9).Value = 7186.5
$ws.Cells.Item(40, 11).Value = 7186.5
$ws.Cells.Item(40, 13).Value = -7050.5
$ws.Cells.Item(133, 8).Value = 81592
$ws.Cells.Item(133, 9).Value = 0
$ws.Cells.Item(133, 10).Value = 81592
$ws.Cells.Item(133, 11).Value = 0
$ws.Cells.Item(133, 12).Value = $null
$ws.Cells.Item(133, 13).Value = 81592
$ws.Cells.Item(133, 14).Value = -86652
$ws.Cells.Item(136, 8).Value = 20836414
$ws.Cells.Item(136, 9).Value = 12503398
$ws.Cells.Item(136, 10).Value = 62501496
$ws.Cells.Item(136, 11).Value = 37510194
$ws.Cells.Item(136, 12).Value = 187504488
$ws.Cells.Item(136, 13).Value = -37507644
$ws.Cells.Item(136, 14).Value = -187509588
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 1400.6875
$ws.Cells.Item(107, 9).Value = 1288.25
$ws.Cells.Item(107, 11).Value = 3864.75
$ws.Cells.Item(107, 13).Value = -1944.75
$ws.Cells.Item(122, 8).Value = 50358.305
$ws.Cells.Item(122, 9).Value = 1361.7646
$ws.Cells.Item(122, 10).Value = 189181.83
$ws.Cells.Item(122, 11).Value = 4085.2938
$ws.Cells.Item(122, 12).Value = 567545.49
$ws.Cells.Item(122, 13).Value = -1635.2938
$ws.Cells.Item(122, 14).Value = -572445.49
$ws.Cells.Item(132, 8).Value = 4763831
$ws.Cells.Item(132, 9).Value = 5954077.5
$ws.Cells.Item(132, 10).Value = 2845.5715
$ws.Cells.Item(132, 11).Value = 17862232.5
$ws.Cells.Item(132, 12).Value = 8536.7145
$ws.Cells.Item(132, 13).Value = -17859702.5
$ws.Cells.Item(132, 14).Value = -13596.7145
